# Updates the "cryptos" price/volume table with the latest scraped figures.
# Many of the Price/Volume(1h) cells hold numeric-looking text (e.g. "289.35",
# "1.143") that must stay plain text (as in the source data), so each write
# temporarily forces a Text number format, assigns the value, then clears the
# format again so the cell's style matches the original (unstyled) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '22.551.65'
Set-TextValue 'E2' '  +0.30%  '

Set-TextValue 'D3' '1.576.85'
Set-TextValue 'E3' '  +0.35%  '

Set-TextValue 'E4' '  +0.00%  '

Set-TextValue 'E5' '  -0.02%  '

Set-TextValue 'D6' '289.35'
Set-TextValue 'E6' '  -0.61%  '

Set-TextValue 'D7' '0.3713'
Set-TextValue 'E7' '  +0.03%  '

Set-TextValue 'D8' '48.58'
Set-TextValue 'E8' '  -2.79%  '

Set-TextValue 'D9' '0.3351'
Set-TextValue 'E9' '  -0.73%  '

Set-TextValue 'D10' '1.143'
Set-TextValue 'E10' '  -0.33%  '

Set-TextValue 'D11' '0.07499'
Set-TextValue 'E11' '  -0.51%  '

Set-TextValue 'E12' '  +0.03%  '

Set-TextValue 'D13' '21.02'
Set-TextValue 'E13' '  -0.58%  '

Set-TextValue 'D14' '5.998'
Set-TextValue 'E14' '  -0.22%  '

Set-TextValue 'D15' '6.954'
Set-TextValue 'E15' '  -0.06%  '

Set-TextValue 'D16' '1.582.65'
Set-TextValue 'E16' '  +0.78%  '

Set-TextValue 'D17' '0.00001122'
Set-TextValue 'E17' '  +0.33%  '

Set-TextValue 'D18' '88.59'
Set-TextValue 'E18' '  -2.08%  '

Set-TextValue 'D19' '0.06775'
Set-TextValue 'E19' '  -0.14%  '

Set-TextValue 'B20' 'Dai'
Set-TextValue 'C20' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D20' '1.002'
Set-TextValue 'E20' '  +0.00%  '

Set-TextValue 'B21' 'Uniswap'
Set-TextValue 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D21' '6.412'
Set-TextValue 'E21' '  +1.14%  '

Set-TextValue 'D22' '16.56'
Set-TextValue 'E22' '  +0.93%  '

Set-TextValue 'D23' '12.16'
Set-TextValue 'E23' '  -0.44%  '

Set-TextValue 'D24' '22.551.22'
Set-TextValue 'E24' '  +0.35%  '

Set-TextValue 'E25' '  +1.31%  '

Set-TextValue 'D26' '2.595'
Set-TextValue 'E26' '  -0.87%  '

Set-TextValue 'D27' '152.79'
Set-TextValue 'E27' '  +2.52%  '

Set-TextValue 'E28' '  -1.31%  '

Set-TextValue 'D29' '5.014'
Set-TextValue 'E29' '  -1.14%  '

Set-TextValue 'D30' '124.44'
Set-TextValue 'E30' '  -0.49%  '

Set-TextValue 'D31' '1.759.65'
Set-TextValue 'E31' '  +0.66%  '

Set-TextValue 'D32' '1.058'
Set-TextValue 'E32' '  -1.03%  '

Set-TextValue 'D33' '6.188'
Set-TextValue 'E33' '  -0.24%  '

Set-TextValue 'D34' '2.018'
Set-TextValue 'E34' '  +0.26%  '

Set-TextValue 'D35' '9.717'
Set-TextValue 'E35' '  -0.71%  '

Set-TextValue 'D36' '0.08341'
Set-TextValue 'E36' '  -0.06%  '

Set-TextValue 'D37' '0.02466'
Set-TextValue 'E37' '  -0.42%  '

Set-TextValue 'D38' '0.2294'
Set-TextValue 'E38' '  -0.30%  '

Set-TextValue 'D39' '5.427'
Set-TextValue 'E39' '  -0.12%  '

Set-TextValue 'D40' '0.06399'
Set-TextValue 'E40' '  -2.17%  '

Set-TextValue 'D41' '1.301'
Set-TextValue 'E41' '  -4.70%  '

Set-TextValue 'D42' '11.38'
Set-TextValue 'E42' '  +0.86%  '

Set-TextValue 'D43' '0.6348'
Set-TextValue 'E43' '  +2.29%  '

Set-TextValue 'D44' '1.001'
Set-TextValue 'E44' '  +0.01%  '

Set-TextValue 'D45' '14.02'
Set-TextValue 'E45' '  -0.52%  '

Set-TextValue 'D46' '0.6191'
Set-TextValue 'E46' '  +6.02%  '

Set-TextValue 'D47' '3.789'
Set-TextValue 'E47' '  -0.49%  '

Set-TextValue 'D48' '2.065'
Set-TextValue 'E48' '  -0.21%  '

Set-TextValue 'D49' '125.24'
Set-TextValue 'E49' '  -3.06%  '

Set-TextValue 'D50' '1.224'
Set-TextValue 'E50' '  -0.19%  '

Set-TextValue 'D51' '0.07283'
Set-TextValue 'E51' '  -0.33%  '
